# Auto-generated edit script: update crypto price/volume data
# per commit "Updated cryptos list on Fri Feb  2 09:45:01 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.005.08"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.311.62"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.58"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.99"
$ws.Range("E6").Value = "  +5.68%  "
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  +4.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.18"
$ws.Range("E10").Value = "  +5.52%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.92"
$ws.Range("E13").Value = "  +15.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.687.25"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.302.79"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.950.00"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  +8.35%  "
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.85"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.40"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("E24").Value = "  +13.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.72"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.08"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.05"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.24"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.08"
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.003.35"
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("E43").Value = "  -6.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0288"
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("E45").Value = "  +7.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.55"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.42"
$ws.Range("E48").Value = "  +6.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.530.55"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").Value = "  +4.94%  "
$ws.Range("E51").Value = "  +0.70%  "
